# Update "想去人数" (interested-count) values on the 展览 (sheet 1),
# 演出 (sheet 2) and 全部类型 (sheet 4) worksheets to reflect newly
# scraped totals.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(3, 6).Value  = 13566   # F3:  13552 -> 13566
$ws1.Cells.Item(6, 6).Value  = 25      # F6:  21    -> 25
$ws1.Cells.Item(8, 6).Value  = 146     # F8:  143   -> 146
$ws1.Cells.Item(9, 6).Value  = 123     # F9:  121   -> 123
$ws1.Cells.Item(10, 6).Value = 82      # F10: 81    -> 82
$ws1.Cells.Item(13, 6).Value = 13568   # F13: 13564 -> 13568
$ws1.Cells.Item(16, 6).Value = 8968    # F16: 8967  -> 8968
$ws1.Cells.Item(17, 6).Value = 9       # F17: 8     -> 9
$ws1.Cells.Item(18, 6).Value = 8062    # F18: 8059  -> 8062
$ws1.Cells.Item(21, 6).Value = 150     # F21: 149   -> 150
$ws1.Cells.Item(23, 6).Value = 148     # F23: 147   -> 148
$ws1.Cells.Item(31, 6).Value = 190     # F31: 188   -> 190

# --- Sheet 2: 演出 -----------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(2, 6).Value  = 41      # F2:  40    -> 41

# --- Sheet 4: 全部类型 -------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(3, 6).Value  = 13566   # F3:  13552 -> 13566
$ws4.Cells.Item(6, 6).Value  = 25      # F6:  21    -> 25
$ws4.Cells.Item(8, 6).Value  = 146     # F8:  143   -> 146
$ws4.Cells.Item(9, 6).Value  = 123     # F9:  121   -> 123
$ws4.Cells.Item(10, 6).Value = 82      # F10: 81    -> 82
$ws4.Cells.Item(13, 6).Value = 13568   # F13: 13564 -> 13568
$ws4.Cells.Item(16, 6).Value = 8968    # F16: 8967  -> 8968
$ws4.Cells.Item(17, 6).Value = 9       # F17: 8     -> 9
$ws4.Cells.Item(18, 6).Value = 8062    # F18: 8059  -> 8062
$ws4.Cells.Item(21, 6).Value = 150     # F21: 149   -> 150
$ws4.Cells.Item(23, 6).Value = 148     # F23: 147   -> 148
$ws4.Cells.Item(29, 6).Value = 41      # F29: 40    -> 41
$ws4.Cells.Item(33, 6).Value = 190     # F33: 188   -> 190

$wb.Save()
